$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right before the "总计" (totals) sheet.
#    Copying an existing quarter sheet (same 8-column layout) instead of
#    adding a blank one so sheetPr / pageMargins / header+index styles come
#    along for free; its data cells are then overwritten below.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheetPos = $totalSheet.Index
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($totalSheet)
$newSheet = $wb.Worksheets.Item($totalSheetPos)
$newSheet.Name = "2022-Q1"

# Forces a cell to be stored as text (so fund codes like "002076" or
# figures like "0.2310" keep their leading/trailing zeros instead of being
# auto-coerced to numbers), then drops the temporary "@" number-format
# back to the default style so no stray formatting is left on the cell.
function Set-TextCell($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.Style = "Normal"
}

# ---- header row (labels change, formatting already copied) ----
$newSheet.Cells.Item(1, 2).Value2 = "基金代码"
$newSheet.Cells.Item(1, 3).Value2 = "基金名称"
$newSheet.Cells.Item(1, 4).Value2 = "基金规模"
$newSheet.Cells.Item(1, 5).Value2 = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value2 = "仓位占比"
$newSheet.Cells.Item(1, 7).Value2 = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value2 = "仓位排名"

# ---- row 2: 002076 浙商中证500指数增强A ----
$newSheet.Cells.Item(2, 1).Value2 = 0
Set-TextCell $newSheet 2 2 "002076"
Set-TextCell $newSheet 2 3 "浙商中证500指数增强A"
Set-TextCell $newSheet 2 4 "14.53"
Set-TextCell $newSheet 2 5 "93.68"
Set-TextCell $newSheet 2 6 "1.59"
Set-TextCell $newSheet 2 7 "0.2310"
$newSheet.Cells.Item(2, 8).Value2 = 4

# ---- row 3: 007386 浙商中证500指数增强C ----
$newSheet.Cells.Item(3, 1).Value2 = 1
Set-TextCell $newSheet 3 2 "007386"
Set-TextCell $newSheet 3 3 "浙商中证500指数增强C"
Set-TextCell $newSheet 3 4 "3.38"
Set-TextCell $newSheet 3 5 "93.68"
Set-TextCell $newSheet 3 6 "1.59"
Set-TextCell $newSheet 3 7 "0.0537"
$newSheet.Cells.Item(3, 8).Value2 = 4

# ---------------------------------------------------------------------------
# 2) Update the "总计" sheet: insert a new top data row for 2022-Q1 and
#    renumber the existing index column.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert(-4121, 0)

# restyle the new index cell (A2) to match the rest of the index column
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(2, 1).PasteSpecial(-4122)

$ws.Cells.Item(2, 1).Value2 = 0
$ws.Cells.Item(2, 2).Value2 = "2022-Q1"
$ws.Cells.Item(2, 2).Style = "Normal"
$ws.Cells.Item(2, 3).Value2 = 2
$ws.Cells.Item(2, 3).Style = "Normal"
$ws.Cells.Item(2, 4).Value2 = 0.28
$ws.Cells.Item(2, 4).Style = "Normal"

# renumber the index column for the rows that shifted down
$ws.Cells.Item(3, 1).Value2 = 1
$ws.Cells.Item(4, 1).Value2 = 2

# Keep the original active sheet selection ("2020-Q4") intact.
$wb.Worksheets.Item("2020-Q4").Activate()
